$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing content (keep column widths / sheet-level props)
$ws.Cells.Clear()

# Row 1: section title
$ws.Range("A1").Value = "Bad Drivers"
$ws.Range("A1").Font.Bold = $true

# Row 2: Bad Drivers column headers
$ws.Range("A2").Value = "Adapter-Driver"
$ws.Range("B2").Value = "Client Count"
$ws.Range("C2").Value = "Critical Minutes"
$ws.Range("D2").Value = "Good Roaming Calculation (%)"
$ws.Range("A2:D2").Borders.Item(9).LineStyle = 1
$ws.Range("B2:D2").HorizontalAlignment = -4152

# Bad Drivers data rows
$ws.Range("A3").Value = 'Intel(R) Dual Band Wireless-AC 8265 - 20.70.27.1'
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 788
$ws.Range("D3").Value = 94.2
$ws.Range("A4").Value = 'Intel(R) Dual Band Wireless-AC 8265 - 20.70.30.1'
$ws.Range("B4").Value = 25
$ws.Range("C4").Value = 2659
$ws.Range("D4").Value = 96.6
$ws.Range("A5").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 23.80.0.7'
$ws.Range("B5").Value = 28
$ws.Range("C5").Value = 2533
$ws.Range("D5").Value = 97.8
$ws.Range("A6").Value = 'Intel(R) Wi-Fi 6 AX201 160MHz - 23.40.0.4'
$ws.Range("B6").Value = 6
$ws.Range("C6").Value = 413
$ws.Range("D6").Value = 98.2
$ws.Range("A7").Value = 'Intel(R) Wi-Fi 6 AX201 160MHz - 23.80.0.7'
$ws.Range("B7").Value = 55
$ws.Range("C7").Value = 2788
$ws.Range("D7").Value = 98.6
$ws.Range("A8").Value = 'Intel(R) Wi-Fi 6E AX211 160MHz - 22.220.0.4'
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 95
$ws.Range("D8").Value = 98.6
$ws.Range("A9").Value = 'Intel(R) Dual Band Wireless-AC 8265 - 20.70.18.2'
$ws.Range("B9").Value = 21
$ws.Range("C9").Value = 1732
$ws.Range("D9").Value = 98.9
$ws.Range("B3:D9").HorizontalAlignment = -4152

# Totals row
$ws.Range("A10").Value = "Totals:"
$ws.Range("A10").Font.Bold = $true
$ws.Range("B10").Value = 137
$ws.Range("C10").Value = 11008
$ws.Range("B10:C10").Font.Bold = $true
$ws.Range("B10:C10").NumberFormat = "#,##0"

# Good Drivers section title
$ws.Range("A16").Value = "Good Drivers (Roaming > 99.8%)"
$ws.Range("A16").Font.Bold = $true

# Good Drivers column headers
$ws.Range("A17").Value = "Adapter-Driver"
$ws.Range("B17").Value = "Total Samples"
$ws.Range("C17").Value = ""
$ws.Range("D17").Value = "Good Roaming Calculation (%)"
$ws.Range("E17").Value = "Driver Vintage"
$ws.Range("A17:E17").Borders.Item(9).LineStyle = 1
$ws.Range("B17").HorizontalAlignment = -4152
$ws.Range("D17:E17").HorizontalAlignment = -4152

# Good Drivers data rows
$ws.Range("A18").Value = 'Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3'
$ws.Range("B18").Value = 11128
$ws.Range("C18").Value = ""
$ws.Range("D18").Value = 100
$ws.Range("E18").ClearContents()
$ws.Range("A19").Value = 'Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4'
$ws.Range("B19").Value = 486214
$ws.Range("C19").Value = ""
$ws.Range("D19").Value = 100
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '2024-11-10'
$ws.Range("A20").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 23.70.2.3'
$ws.Range("B20").Value = 18721
$ws.Range("C20").Value = ""
$ws.Range("D20").Value = 99.9
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '2024-07-23'
$ws.Range("A21").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.250.10.1'
$ws.Range("B21").Value = 69578
$ws.Range("C21").Value = ""
$ws.Range("D21").Value = 99.9
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '2023-08-14'
$ws.Range("A22").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.230.0.8'
$ws.Range("B22").Value = 338880
$ws.Range("C22").Value = ""
$ws.Range("D22").Value = 99.9
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '2023-05-08'
$ws.Range("A23").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.200.0.6'
$ws.Range("B23").Value = 143869
$ws.Range("C23").Value = ""
$ws.Range("D23").Value = 99.9
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '2023-01-16'
$ws.Range("A24").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.190.0.4'
$ws.Range("B24").Value = 287148
$ws.Range("C24").Value = ""
$ws.Range("D24").Value = 99.9
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '2022-11-22'
$ws.Range("A25").Value = 'Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1'
$ws.Range("B25").Value = 11140
$ws.Range("C25").Value = ""
$ws.Range("D25").Value = 99.9
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '2022-08-29'
$ws.Range("A26").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.160.0.4'
$ws.Range("B26").Value = 96526
$ws.Range("C26").Value = ""
$ws.Range("D26").Value = 99.9
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '2022-08-13'
$ws.Range("A27").Value = 'Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3'
$ws.Range("B27").Value = 14487
$ws.Range("C27").Value = ""
$ws.Range("D27").Value = 99.9
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '2022-05-23'
$ws.Range("A28").Value = 'Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1'
$ws.Range("B28").Value = 265400
$ws.Range("C28").Value = ""
$ws.Range("D28").Value = 99.9
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '2022-05-01'
$ws.Range("A29").Value = 'Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9'
$ws.Range("B29").Value = 79953
$ws.Range("C29").Value = ""
$ws.Range("D29").Value = 99.9
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '2021-08-18'
$ws.Range("A30").Value = 'Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1'
$ws.Range("B30").Value = 35355
$ws.Range("C30").Value = ""
$ws.Range("D30").Value = 99.9
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '2021-04-27'
$ws.Range("A31").Value = 'Intel(R) Dual Band Wireless-AC 8265 - 22.30.0.11'
$ws.Range("B31").Value = 172690
$ws.Range("C31").Value = ""
$ws.Range("D31").Value = 99.9
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '2021-01-19'
$ws.Range("A32").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.30.0.11'
$ws.Range("B32").Value = 67111
$ws.Range("C32").Value = ""
$ws.Range("D32").Value = 100
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '2021-01-19'
$ws.Range("A33").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.10.0.7'
$ws.Range("B33").Value = 68450
$ws.Range("C33").Value = ""
$ws.Range("D33").Value = 100
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '2020-10-19'
$ws.Range("A34").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.0.1.1'
$ws.Range("B34").Value = 15734
$ws.Range("C34").Value = ""
$ws.Range("D34").Value = 99.9
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '2020-09-28'
$ws.Range("A35").Value = 'Intel(R) Dual Band Wireless-AC 8265 - 22.0.1.1'
$ws.Range("B35").Value = 52096
$ws.Range("C35").Value = ""
$ws.Range("D35").Value = 100
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '2020-09-28'
$ws.Range("A36").Value = 'Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2'
$ws.Range("B36").Value = 65425
$ws.Range("C36").Value = ""
$ws.Range("D36").Value = 100
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '2020-08-05'
$ws.Range("A37").Value = 'Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6'
$ws.Range("B37").Value = 117653
$ws.Range("C37").Value = ""
$ws.Range("D37").Value = 100
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '2020-01-06'
$ws.Range("A38").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 21.60.2.1'
$ws.Range("B38").Value = 26241
$ws.Range("C38").Value = ""
$ws.Range("D38").Value = 100
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '2019-12-14'
$ws.Range("A39").Value = 'Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1'
$ws.Range("B39").Value = 56018
$ws.Range("C39").Value = ""
$ws.Range("D39").Value = 100
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '2019-12-14'
$ws.Range("A40").Value = 'Intel(R) Dual Band Wireless-AC 8265 - 20.70.11.3'
$ws.Range("B40").Value = 161874
$ws.Range("C40").Value = ""
$ws.Range("D40").Value = 100
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '2019-09-05'
$ws.Range("A41").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 21.40.2.2'
$ws.Range("B41").Value = 90508
$ws.Range("C41").Value = ""
$ws.Range("D41").Value = 99.9
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '2019-08-31'
$ws.Range("A42").Value = 'Intel(R) Dual Band Wireless-AC 8265 - 20.70.12.5'
$ws.Range("B42").Value = 154175
$ws.Range("C42").Value = ""
$ws.Range("D42").Value = 99.9
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '2019-08-25'
$ws.Range("A43").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 21.30.4.1'
$ws.Range("B43").Value = 13016
$ws.Range("C43").Value = ""
$ws.Range("D43").Value = 100
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '2019-07-29'
$ws.Range("A44").Value = 'Intel(R) Dual Band Wireless-AC 8265 - 20.70.10.2'
$ws.Range("B44").Value = 20227
$ws.Range("C44").Value = ""
$ws.Range("D44").Value = 100
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '2019-05-11'
$ws.Range("A45").Value = 'Intel(R) Dual Band Wireless-AC 8265 - 20.70.9.1'
$ws.Range("B45").Value = 34065
$ws.Range("C45").Value = ""
$ws.Range("D45").Value = 100
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '2019-04-28'
$ws.Range("A46").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 21.10.1.2'
$ws.Range("B46").Value = 52515
$ws.Range("C46").Value = ""
$ws.Range("D46").Value = 100
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '2019-04-23'
$ws.Range("A47").Value = 'Intel(R) Dual Band Wireless-AC 8265 - 20.70.8.1'
$ws.Range("B47").Value = 48540
$ws.Range("C47").Value = ""
$ws.Range("D47").Value = 100
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '2019-03-16'
$ws.Range("A48").Value = 'Intel(R) Dual Band Wireless-AC 8265 - 20.70.5.2'
$ws.Range("B48").Value = 184564
$ws.Range("C48").Value = ""
$ws.Range("D48").Value = 99.9
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '2018-11-25'
$ws.Range("A49").Value = 'Intel(R) Dual Band Wireless-AC 8265 - 20.50.0.4'
$ws.Range("B49").Value = 14221
$ws.Range("C49").Value = ""
$ws.Range("D49").Value = 100
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '2018-05-08'
$ws.Range("A50").Value = 'Intel(R) Dual Band Wireless-AC 8265 - 20.30.1.2'
$ws.Range("B50").Value = 23765
$ws.Range("C50").Value = ""
$ws.Range("D50").Value = 100
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '2018-01-09'
$ws.Range("B18:B50").NumberFormat = "#,##0"
$ws.Range("B18:B50").HorizontalAlignment = -4152
$ws.Range("D18:E50").HorizontalAlignment = -4152

